$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Helper: insert a blank row at row index $r, copy the formatting of row
# $templateRow into it (so the new row's cell styles match the rest of the
# table), and return nothing. The workbook's data rows all use a single
# uniform style, so any existing data row can serve as the template.
# ---------------------------------------------------------------------------
function Insert-FormattedRow($r, $templateRow) {
    $ws.Rows.Item($r).Insert()
    $ws.Range("A$templateRow`:K$templateRow").Copy()
    $ws.Range("A$r`:K$r").PasteSpecial(-4122)  # xlPasteFormats
}

# 1) Two new rows after row 8: SpO2 observation, Breathing rate observation.
Insert-FormattedRow 9 10
Insert-FormattedRow 10 11

$ws.Cells.Item(9,1).Value = "tr-br-breathing-finding-oxygen-saturation-observation"
$ws.Cells.Item(9,2).Value = "Observation of SpO2"
$ws.Cells.Item(9,5).Value = "SNOMED CT#103228002"
$ws.Cells.Item(9,7).Value = "dateTime" + [char]0x135 + ", Period" + [char]0x135 + ", Timing" + [char]0x135 + ", instant" + [char]0x135
$ws.Cells.Item(9,8).Value = "Range" + [char]0x135
$ws.Cells.Item(9,9).Value = "optional"

$ws.Cells.Item(10,1).Value = "tr-br-breathing-finding-respiration-rate-observation"
$ws.Cells.Item(10,2).Value = "Observation of Breathing rate"
$ws.Cells.Item(10,5).Value = "SNOMED CT#86290005"
$ws.Cells.Item(10,7).Value = "dateTime" + [char]0x135 + ", Period" + [char]0x135 + ", Timing" + [char]0x135 + ", instant" + [char]0x135
$ws.Cells.Item(10,8).Value = "CodeableConcept" + [char]0x135
$ws.Cells.Item(10,9).Value = "optional"

# 2) Two new rows before the (old) row 13 "tr-c-skin-assessment-observation":
#    Cardiac Arrest block. After the first insertion above, the old row 13
#    now lives at row 15, so we insert immediately above that.
Insert-FormattedRow 15 14
Insert-FormattedRow 16 14

$ws.Cells.Item(15,1).Value = "tr-circulation-cardiac-arrest-observation"
$ws.Cells.Item(15,2).Value = "Cardiac Arrest"
$ws.Cells.Item(15,5).Value = "SNOMED CT#410429000"
$ws.Cells.Item(15,7).Value = "dateTime" + [char]0x135 + ", Period" + [char]0x135 + ", Timing" + [char]0x135 + ", instant" + [char]0x135
$ws.Cells.Item(15,8).Value = "Quantity, CodeableConcept, string, boolean, integer, Range, Ratio, SampledData, time, dateTime, Period"
$ws.Cells.Item(15,9).Value = "optional"

$ws.Cells.Item(16,2).Value = "Cardiac Arrest"
$ws.Cells.Item(16,5).Value = "SNOMED CT#45007003"
$ws.Cells.Item(16,8).Value = "CodeableConcept"
$ws.Cells.Item(16,9).Value = "optional"

# 3) The old row 13 (now row 17) is renamed from tr-c-skin-assessment-observation
#    to tr-circulation-skin-assessment-observation.
$ws.Cells.Item(17,1).Value = "tr-circulation-skin-assessment-observation"

# 4) Two new rows appended to the skin-assessment block, after the old row 16
#    (LOINC#39106-0), which now lives at row 20.
Insert-FormattedRow 21 20
Insert-FormattedRow 22 20

$ws.Cells.Item(21,2).Value = "Skin assessment"
$ws.Cells.Item(21,5).Value = "SNOMED CT#409055009"
$ws.Cells.Item(21,8).Value = "boolean"
$ws.Cells.Item(21,9).Value = "optional"

$ws.Cells.Item(22,2).Value = "Skin assessment"
$ws.Cells.Item(22,5).Value = "SNOMED CT#50960005"
$ws.Cells.Item(22,8).Value = "CodeableConcept"
$ws.Cells.Item(22,9).Value = "optional"

Write-Host "done"
